$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously only had "MONTEREY AREA TOTALS" in column B.
# Now column A gets "MONTEREY AREA TOTALS" and column B gets "Totals".
$ws.Range("A2").Value = "MONTEREY AREA TOTALS"
$ws.Range("B2").Value = "Totals"

# Update the active selection to B3, matching the saved view state.
$ws.Range("B3").Select()
